# fix(import): add new features to import
#
# Inserts a new "Quelle est la situation professionelle de la personne ?"
# column in the usagers-import-test template (Feuil1), right before the
# existing column U, with a couple of sample answers on the demo rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- 1. Insert a new column before U -----------------------------------
# Everything that used to live in U:AX (headers, data, named ranges that
# reference it) shifts right by one, to V:AY.
$ws.Columns("U").Insert()

# --- 2. New header cell (U1) --------------------------------------------
$header = $ws.Range("U1")
$header.Value = "Quelle est la situation professionelle de la personne ?"
$header.Font.Name = "Arial"
$header.Font.Bold = $true
$header.Font.Size = 12.5
$header.Font.Color = 2500134      # RGB(0x26,0x26,0x26)
$header.Interior.PatternColor = 16777215   # RGB(255,255,255)
$header.Interior.Color = 16777215
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108        # xlCenter
$header.VerticalAlignment = -4108          # xlCenter
$header.WrapText = $true
$header.NumberFormat = "@"

# --- 3. Sample answers on the first few demo rows ------------------------
$data = $ws.Range("U2:U4")
$data.Font.Name = "Arial"
$data.Font.Size = 13
$data.Interior.PatternColor = 15987699     # RGB(0xF3,0xF3,0xF3)
$data.Interior.Color = 15987699
$data.Borders.LineStyle = 1
$data.Borders.Weight = 2
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108
$data.NumberFormat = "@"

$ws.Range("U2").Value = "SALARIE"
$ws.Range("U3").Value = "SALARIE"
$ws.Range("U4").Value = "FRANCE_TRAVAIL"

# --- 4. Row heights for the rows that now carry the taller answer text ---
$ws.Rows("2:4").RowHeight = 17

# --- 5. Keep the hidden _FilterDatabase name in sync with the new column -
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Feuil1!_FilterDatabase") {
        $n.RefersTo = "=Feuil1!`$A`$1:`$AY`$1"
    }
}

# --- 6. Refresh the view: scroll over to the new column & select a cell -
$ws.Activate()
$ws.Range("U9").Select()

Write-Host "Added situation professionnelle column"
